$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Change 1: cell (row 18, col 3) "29.10" -> split into two runs "29." + "09" ---
$cell1 = $t.Cell(18, 3)
$range1 = $cell1.Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="00BA5800" w:rsidRPr="0054465A" w:rsidRDefault="003245EB" w:rsidP="007A5C54">' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                  '<w:sz w:val="28"/>' +
                  '<w:szCs w:val="28"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                  '<w:sz w:val="28"/>' +
                  '<w:szCs w:val="28"/>' +
                '</w:rPr>' +
                '<w:t>29.</w:t>' +
              '</w:r>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                  '<w:sz w:val="28"/>' +
                  '<w:szCs w:val="28"/>' +
                '</w:rPr>' +
                '<w:t>09</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$range1.InsertXML($xml1)

# --- Change 2: cell (row 20, col 3) empty paragraph -> add run with "03.10" ---
$cell2 = $t.Cell(20, 3)
$range2 = $cell2.Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="00BA5800" w:rsidRPr="0054465A" w:rsidRDefault="00BA5800" w:rsidP="00C12D0E">' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                  '<w:sz w:val="28"/>' +
                  '<w:szCs w:val="28"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                  '<w:sz w:val="28"/>' +
                  '<w:szCs w:val="28"/>' +
                '</w:rPr>' +
                '<w:t>03.10</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$range2.InsertXML($xml2)

Write-Host "Edits applied."
